# Insert a new data row at row 649 of Sheet1, pushing existing rows
# 649:702 down to 650:703, and populate the newly inserted row with its
# own values (date, quality, volume, prices, origin, price/kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 649 (this shifts rows
# 649-702 down to 650-703 automatically, carrying their formatting).
$ws.Rows.Item(649).Insert()

# Populate the newly inserted row 649 with the new record's data.
$ws.Cells.Item(649, 1).Value2  = 4
$ws.Cells.Item(649, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(649, 3).Value2  = "Los Lagos"
$ws.Cells.Item(649, 4).Value2  = 45106
$ws.Cells.Item(649, 5).Value2  = 10
$ws.Cells.Item(649, 6).Value2  = 100114001
$ws.Cells.Item(649, 7).Value2  = "Papa"
$ws.Cells.Item(649, 8).Value2  = "Patagonia"
$ws.Cells.Item(649, 9).Value2  = "1a (guarda)"
$ws.Cells.Item(649, 10).Value2 = 150
$ws.Cells.Item(649, 11).Value2 = 16000
$ws.Cells.Item(649, 12).Value2 = 16000
$ws.Cells.Item(649, 13).Value2 = 16000
$ws.Cells.Item(649, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(649, 15).Value2 = "Provincia de Llanquihue"
$ws.Cells.Item(649, 16).Value2 = 640
$ws.Cells.Item(649, 17).Value2 = 25
$ws.Cells.Item(649, 18).Value2 = "Hortaliza"
